# Update column L (correct_ans) abbreviations to full words
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$correctAns = @(
    "left","center","right","center","right","left","left","center","right","right",
    "left","center","left","right","center","left","right","center","center","right",
    "left","center","left","right","center","right","left","center","left","right",
    "right","left","center","right","center","left","left","center","right","center",
    "right","left","left","right","center","center","right","left","right","center",
    "left","right","left","center","left","right","center","left","center","right",
    "left","center","right","right","left","center","center","left","right","left",
    "right","center","right","center","left","center","left","right","left","right",
    "center","center","right","left","center","right","left","left","right","center",
    "right","center","left","left","right","center","right","center","left","center",
    "left","right","left","right","center","right","center","left","right","center",
    "left","right","center","left","right","center","left","center","right","left",
    "right","center","left","right","center","left","right","left","center","center",
    "left","right","right","center","left","left","right","center","center","left",
    "right","left","right","center","right","center","left","right","left","center",
    "right","center","left","right","left","center","center","left","right","left",
    "right","center","left","center","right","left","center","right","left","right",
    "center","center","left","right","center","left","right","left","right","center",
    "right","center","left","right","left","center","center","left","right","right",
    "center","left","left","right","center","center","left","right","right","center",
    "left","left","center","right","center","left","right","center","left","right",
    "right","left","center","left","center","right","left","right","center","left",
    "center","right","left","center","right","left","center","right","right","center",
    "left","right","left","center","center","left","right","left","center","right",
    "left","right","center","right","center","left","center","left","right","center",
    "left","right","left","right","center","right","left","center","center","right",
    "left","center","left","right","right","center","left","center","right","left",
    "right","center","left","right","left","center","left","right","center","right",
    "center","left","center","right","left","right","center","left","right","left",
    "center","left","center","right","left","center","right","center","right","left",
    "right","left","center","right","left","center","right","left","center","center",
    "right","left","right","left","center","center","right","left","right","left",
    "center","right","center","left","center","right","left","right","center","left",
    "center","left","right","center","right","left","center","left","right","center",
    "right","left","left","right","center","right","left","center","left","right",
    "center","right","left","center","right","left","center","center","right","left"
)

$startRow = 2
for ($i = 0; $i -lt $correctAns.Length; $i++) {
    $ws.Cells.Item($startRow + $i, 12).Value = $correctAns[$i]
}

# Rename "face" stimulus category to "book" across all referencing cells
$faceToBook = @(
    @{Cell="D17"; Value="book//book_18.jpg"},
    @{Cell="D82"; Value="book//book_19.jpg"},
    @{Cell="D94"; Value="book//book_16.jpg"},
    @{Cell="D111"; Value="book//book_04.jpg"},
    @{Cell="D116"; Value="book//book_32.jpg"},
    @{Cell="D117"; Value="book//book_33.jpg"},
    @{Cell="D118"; Value="book//book_07.jpg"},
    @{Cell="D133"; Value="book//book_12.jpg"},
    @{Cell="D145"; Value="book//book_16.jpg"},
    @{Cell="D161"; Value="book//book_35.jpg"},
    @{Cell="D169"; Value="book//book_24.jpg"},
    @{Cell="D179"; Value="book//book_01.jpg"},
    @{Cell="D205"; Value="book//book_28.jpg"},
    @{Cell="D211"; Value="book//book_15.jpg"},
    @{Cell="D213"; Value="book//book_06.jpg"},
    @{Cell="D227"; Value="book//book_25.jpg"},
    @{Cell="D229"; Value="book//book_17.jpg"},
    @{Cell="D237"; Value="book//book_05.jpg"},
    @{Cell="D239"; Value="book//book_19.jpg"},
    @{Cell="D243"; Value="book//book_18.jpg"},
    @{Cell="D254"; Value="book//book_13.jpg"},
    @{Cell="D271"; Value="book//book_11.jpg"},
    @{Cell="D276"; Value="book//book_40.jpg"},
    @{Cell="D296"; Value="book//book_30.jpg"},
    @{Cell="A302"; Value="book//book_36.jpg"},
    @{Cell="C304"; Value="book//book_36.jpg"},
    @{Cell="C306"; Value="book//book_36.jpg"},
    @{Cell="A308"; Value="book//book_22.jpg"},
    @{Cell="C311"; Value="book//book_22.jpg"},
    @{Cell="A314"; Value="book//book_10.jpg"},
    @{Cell="A320"; Value="book//book_33.jpg"},
    @{Cell="C325"; Value="book//book_33.jpg"},
    @{Cell="A326"; Value="book//book_17.jpg"},
    @{Cell="C330"; Value="book//book_17.jpg"},
    @{Cell="A332"; Value="book//book_20.jpg"},
    @{Cell="A338"; Value="book//book_31.jpg"},
    @{Cell="C340"; Value="book//book_31.jpg"},
    @{Cell="C341"; Value="book//book_31.jpg"},
    @{Cell="A344"; Value="book//book_34.jpg"},
    @{Cell="C349"; Value="book//book_34.jpg"},
    @{Cell="A350"; Value="book//book_13.jpg"},
    @{Cell="C352"; Value="book//book_13.jpg"},
    @{Cell="C353"; Value="book//book_13.jpg"},
    @{Cell="C355"; Value="book//book_13.jpg"},
    @{Cell="A356"; Value="book//book_30.jpg"},
    @{Cell="C360"; Value="book//book_30.jpg"}
)

foreach ($item in $faceToBook) {
    $ws.Range($item.Cell).Value = $item.Value
}
